$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column I (9th column); this shifts the existing
# "Logistics(%)" column (header + data) from I to J.
$ws.Columns.Item(9).Insert()

# New header cell I1: "Distribution channel code" (bold, matching other headers)
$ws.Cells.Item(1, 9).Value = "Distribution channel code"
$ws.Cells.Item(1, 9).Font.Bold = $true

# New data cells for the inserted column I
$ws.Cells.Item(2, 9).Value = "IN"
$ws.Cells.Item(3, 9).Value = "GO"

# Widen the new column to fit its header text
$ws.Columns.Item(9).ColumnWidth = 21.6

# Update the active selection to match the edited workbook (I3)
$ws.Range("I3").Select()
